# Scheduled runner update: refresh computed profit columns (H:N) in each
# leve-profit sheet from the latest market-board averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4785737.5
$ws.Range("I19").Value = 9091901
$ws.Range("J19").Value = 1111.3334
$ws.Range("K19").Value = 9091901
$ws.Range("L19").Value = 1111.3334
$ws.Range("M19").Value = -9091726
$ws.Range("N19").Value = -1461.3334

$ws.Range("H96").Value = 291.04
$ws.Range("I96").Value = 167.07143
$ws.Range("J96").Value = 448.81818
$ws.Range("K96").Value = 501.21429
$ws.Range("L96").Value = 1346.45454
$ws.Range("M96").Value = 871.78571
$ws.Range("N96").Value = -4092.45454

$ws.Range("H100").Value = 1382.1177
$ws.Range("I100").Value = 1200.1111
$ws.Range("J100").Value = 1586.875
$ws.Range("K100").Value = 1200.1111
$ws.Range("L100").Value = 1586.875
$ws.Range("M100").Value = -659.1111000000001
$ws.Range("N100").Value = -2668.875

$ws.Range("H112").Value = 1781
$ws.Range("J112").Value = 1838.7059
$ws.Range("L112").Value = 5516.1177
$ws.Range("N112").Value = -7732.1177

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 2199.4285
$ws.Range("I132").Value = 1932.2264
$ws.Range("J132").Value = 3032.4707
$ws.Range("K132").Value = 5796.6792
$ws.Range("L132").Value = 9097.4121
$ws.Range("M132").Value = -3266.6792
$ws.Range("N132").Value = -14157.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17246026
$ws.Range("I32").Value = 19234338
$ws.Range("K32").Value = 19234338
$ws.Range("M32").Value = -19234051

$ws.Range("H45").Value = 1959.7142
$ws.Range("I45").Value = 1230
$ws.Range("K45").Value = 1230
$ws.Range("M45").Value = -853

$ws.Range("H102").Value = 1143.6316
$ws.Range("I102").Value = 968.5625
$ws.Range("J102").Value = 2077.3333
$ws.Range("K102").Value = 968.5625
$ws.Range("L102").Value = 2077.3333
$ws.Range("M102").Value = 653.4375
$ws.Range("N102").Value = -5321.3333

$ws.Range("H110").Value = 2485.2666
$ws.Range("I110").Value = 1977.9
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 1977.9
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 67.09999999999991
$ws.Range("N110").Value = -7590

$ws.Range("H139").Value = 27980
$ws.Range("J139").Value = 27980
$ws.Range("L139").Value = 27980
$ws.Range("N139").Value = -38260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1921.7727
$ws.Range("I86").Value = 1954.9375
$ws.Range("K86").Value = 1954.9375
$ws.Range("M86").Value = -831.9375

$ws.Range("H89").Value = 1921.7727
$ws.Range("I89").Value = 1954.9375
$ws.Range("K89").Value = 9774.6875
$ws.Range("M89").Value = -4158.6875

$ws.Range("H105").Value = 2630.39
$ws.Range("I105").Value = 1903.2963
$ws.Range("J105").Value = 2899.3152
$ws.Range("K105").Value = 1903.2963
$ws.Range("L105").Value = 2899.3152
$ws.Range("M105").Value = -156.2963
$ws.Range("N105").Value = -6393.3152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H99").Value = 1947.1
$ws.Range("I99").Value = 1941.2222
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1941.2222
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -443.2221999999999
$ws.Range("N99").Value = -4996

$ws.Range("H105").Value = 754.5
$ws.Range("I105").Value = 754.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 754.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 992.5
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 788.6111
$ws.Range("I107").Value = 675.36365
$ws.Range("J107").Value = 966.5714
$ws.Range("K107").Value = 675.36365
$ws.Range("L107").Value = 966.5714
$ws.Range("M107").Value = 1244.63635
$ws.Range("N107").Value = -4806.5714

$ws.Range("H126").Value = 1947.1
$ws.Range("I126").Value = 1941.2222
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5823.6666
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3353.6666
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 379.33334
$ws.Range("J40").Value = 449.16666
$ws.Range("L40").Value = 1796.66664
$ws.Range("N40").Value = -1934.66664

$ws.Range("H131").Value = 827.82855
$ws.Range("I131").Value = 590.63635
$ws.Range("J131").Value = 936.5417
$ws.Range("K131").Value = 1771.90905
$ws.Range("L131").Value = 2809.6251
$ws.Range("M131").Value = 3268.09095
$ws.Range("N131").Value = -12889.6251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3031864
$ws.Range("I126").Value = 4762923
$ws.Range("J126").Value = 2510.25
$ws.Range("K126").Value = 14288769
$ws.Range("L126").Value = 7530.75
$ws.Range("M126").Value = -14286299
$ws.Range("N126").Value = -12470.75

$ws.Range("H132").Value = 3621.5386
$ws.Range("I132").Value = 3655.525
$ws.Range("J132").Value = 3508.25
$ws.Range("K132").Value = 10966.575
$ws.Range("L132").Value = 10524.75
$ws.Range("M132").Value = -8436.575000000001
$ws.Range("N132").Value = -15584.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 63752436
$ws.Range("I100").Value = 2002998
$ws.Range("J100").Value = 166668160
$ws.Range("K100").Value = 2002998
$ws.Range("L100").Value = 166668160
$ws.Range("M100").Value = -2002457
$ws.Range("N100").Value = -166669242

$ws.Range("H122").Value = 4349.3
$ws.Range("I122").Value = 3698.1667
$ws.Range("K122").Value = 11094.5001
$ws.Range("M122").Value = -8644.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2312.6316
$ws.Range("I126").Value = 1998.5
$ws.Range("J126").Value = 3192.2
$ws.Range("K126").Value = 5995.5
$ws.Range("L126").Value = 9576.599999999999
$ws.Range("M126").Value = -3525.5
$ws.Range("N126").Value = -14516.6

$ws.Range("H132").Value = 1665.7646
$ws.Range("I132").Value = 923.74194
$ws.Range("J132").Value = 2815.9
$ws.Range("K132").Value = 2771.22582
$ws.Range("L132").Value = 8447.700000000001
$ws.Range("M132").Value = -241.2258200000001
$ws.Range("N132").Value = -13507.7

$ws.Range("H136").Value = 2833.3394
$ws.Range("I136").Value = 2641.7446
$ws.Range("J136").Value = 3833.889
$ws.Range("K136").Value = 7925.2338
$ws.Range("L136").Value = 11501.667
$ws.Range("M136").Value = -5375.2338
$ws.Range("N136").Value = -16601.667
